$d = $word.ActiveDocument

$para = $d.Paragraphs.Item(1)
$rng = $para.Range
$rng.Text = "timeline"
